$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New shared strings / cell content (order matters for shared-string index assignment)
$ws.Range("J14").Value = "Team"
$ws.Range("B11").Value = "Team: Coloque un mismo número para aquellos estudiantes que hagan grupo en las categorias correspondientes"

# Update selection to match the diff (activeCell G12)
$ws.Range("G12").Select()
